$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Junio de 2020 a las 15:05"

# Apply updated country stats (values refreshed; a handful of countries swapped
# rank/position because their totals changed, so column A is updated too where needed)

# Row 4
$ws.Range("B4").Value = 2424904
$ws.Range("C4").Value = 736
$ws.Range("D4").Value = 1020414
$ws.Range("E4").Value = 1281006
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 123484

# Row 7
$ws.Range("B7").Value = 461444
$ws.Range("C7").Value = 5329
$ws.Range("D7").Value = 261545
$ws.Range("E7").Value = 185330
$ws.Range("G7").Value = 86
$ws.Range("H7").Value = 14569

# Row 14
$ws.Range("D14").Value = 176300
$ws.Range("E14").Value = 7541

# Row 18
$ws.Range("B18").Value = 167267
$ws.Range("C18").Value = 3123
$ws.Range("D18").Value = 112797
$ws.Range("E18").Value = 53083
$ws.Range("G18").Value = 41
$ws.Range("H18").Value = 1387

# Row 31
$ws.Range("B31").Value = 49804
$ws.Range("C31").Value = 82
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 6097

# Row 37
$ws.Range("B37").Value = 40104
$ws.Range("C37").Value = 367
$ws.Range("D37").Value = 26083
$ws.Range("E37").Value = 12478
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 1543

# Row 50
$ws.Range("E50").Value = 5544
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 68

# Row 62
$ws.Range("B62").Value = 13235
$ws.Range("C62").Value = 143
$ws.Range("D62").Value = 12111
$ws.Range("E62").Value = 861

# Row 63
$ws.Range("B63").Value = 12615
$ws.Range("C63").Value = 54
$ws.Range("D63").Value = 11422
$ws.Range("E63").Value = 590

# Row 76
$ws.Range("B76").Value = 6847
$ws.Range("C76").Value = 185
$ws.Range("E76").Value = 2240

# Row 80
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("B80").Value = 5445
$ws.Range("C80").Value = 134
$ws.Range("D80").Value = 2091
$ws.Range("E80").Value = 3095
$ws.Range("G80").Value = 8
$ws.Range("H80").Value = 259

# Row 81
$ws.Range("A81").Value = "Haiti"
$ws.Range("B81").Value = 5324
$ws.Range("C81").Value = 113
$ws.Range("D81").Value = 436
$ws.Range("E81").Value = 4799
$ws.Range("G81").Value = 1
$ws.Range("H81").Value = 89

# Row 87
$ws.Range("B87").Value = 4630
$ws.Range("C87").Value = 13
$ws.Range("D87").Value = 4182
$ws.Range("E87").Value = 396
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 52

# Row 93
$ws.Range("B93").Value = 3676
$ws.Range("C93").Value = 88
$ws.Range("D93").Value = 2297
$ws.Range("E93").Value = 1206
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 173

# Row 101
$ws.Range("A101").Value = "Croacia"
$ws.Range("B101").Value = 2388
$ws.Range("C101").Value = 22
$ws.Range("D101").Value = 2145
$ws.Range("E101").Value = 136
$ws.Range("H101").Value = 107

# Row 102
$ws.Range("A102").Value = "Costa Rica"
$ws.Range("B102").Value = 2368
$ws.Range("D102").Value = 1129
$ws.Range("E102").Value = 1227
$ws.Range("H102").Value = 12

# Row 107
$ws.Range("A107").Value = "Mali"
$ws.Range("B107").Value = 2001
$ws.Range("C107").Value = 23
$ws.Range("D107").Value = 1333
$ws.Range("E107").Value = 556
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 112

# Row 108
$ws.Range("A108").Value = "Sri Lanka"
$ws.Range("B108").Value = 1998
$ws.Range("C108").Value = 7
$ws.Range("D108").Value = 1562
$ws.Range("E108").Value = 425
$ws.Range("H108").Value = 11

# Row 109
$ws.Range("A109").Value = "Estonia"
$ws.Range("B109").Value = 1983
$ws.Range("C109").Value = 1
$ws.Range("D109").Value = 1783
$ws.Range("E109").Value = 131
$ws.Range("H109").Value = 69

# Row 135
$ws.Range("A135").Value = "Benin"
$ws.Range("B135").Value = 902
$ws.Range("C135").Value = 52
$ws.Range("D135").Value = 277
$ws.Range("E135").Value = 612
$ws.Range("H135").Value = 13

# Row 136
$ws.Range("A136").Value = "Uruguay"
$ws.Range("B136").Value = 885
$ws.Range("D136").Value = 815
$ws.Range("E136").Value = 45
$ws.Range("H136").Value = 25

# Row 137
$ws.Range("A137").Value = "Republica del Chad"
$ws.Range("B137").Value = 860
$ws.Range("D137").Value = 757
$ws.Range("E137").Value = 29
$ws.Range("H137").Value = 74

# Row 138
$ws.Range("A138").Value = "Principado de Andorra"
$ws.Range("B138").Value = 855
$ws.Range("D138").Value = 797
$ws.Range("E138").Value = 6
$ws.Range("H138").Value = 52

# Row 202
$ws.Range("A202").Value = "Fiyi"

# Row 203
$ws.Range("A203").Value = "Dominica"

# Row 211
$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Row 212
$ws.Range("A212").Value = "Montserrat"
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1
